{"js": "// Office.js (Word JavaScript API) edit script.\n// Wraps the JSON array of translation-review objects in a ```json fenced\n// code block, fixes the first entry's text, and appends two more entries.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- 1. Insert a \"```json\" paragraph before the current first paragraph ---\nconst firstPara = paragraphs.items[0];\nfirstPara.insertParagraph(\"```json\", Word.InsertLocation.before);\nawait context.sync();\n\n// Reload paragraphs now that a new one has been added at the top.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Paragraph index map (after the ```json insertion):\n// 0: ```json\n// 1: [\n// 2:   {\n// 3:     \"problematic_translated_sentence\": ...\n// 4:     \"respective_source_sentence\": ...\n// 5:     \"errors\": ...\n// 6:   }\n// 7: ]\n\n// --- 2. Replace the text of the three data paragraphs of the first entry ---\nparagraphs.items[3].insertText(\n  '    \"problematic_translated_sentence\": \"\u6240\u6709\u5728\u5f37\u8feb\u4e0b\u505a\u51fa\u7684\u52aa\u529b\u90fd\u9700\u8981\u72a7\u7272\u751f\u547d\u529b\u3002\",',\n  Word.InsertLocation.replace\n);\nparagraphs.items[4].insertText(\n  '    \"respective_source_sentence\": \"Every effort under compulsion demands a sacrifice of life-energy.\",',\n  Word.InsertLocation.replace\n);\nparagraphs.items[5].insertText(\n  '    \"errors\": \"Duplication of \\'\u6240\u6709\\' and \\'\u505a\u51fa\u7684\\' is redundant.\"',\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 3. After the \"errors\" paragraph (index 5), insert the closing \"},\"\n//         of entry 1 and the two new entries, before the closing \"}\" (index 6) ---\nconst closingBrace = paragraphs.items[6]; // the lone \"  }\" paragraph\n\nconst newParaTexts = [\n  '  },',\n  '  {',\n  '    \"problematic_translated_sentence\": \"\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\",',\n  '    \"respective_source_sentence\": \"I never paid such a price.\",',\n  '    \"errors\": \"Mistranslation of \\'paid such a price\\' as \\'\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\\'.\"',\n  '  },',\n  '  {',\n  '    \"problematic_translated_sentence\": \"\u76f8\u53cd\uff0c\u6211\u5f9e\u6211\u7684\u601d\u7dd2\u4e2d\u7372\u76ca\u532a\u6dfa\u3002\",',\n  '    \"respective_source_sentence\": \"On the contrary, I have thrived on my thoughts.\",',\n  '    \"errors\": \"Mistranslation of \\'thrived on my thoughts\\' as \\'\u5f9e\u6211\u7684\u601d\u7dd2\u4e2d\u7372\u76ca\u532a\u6dfa\\'.\"',\n];\n\n// Insert in forward (top-to-bottom) order, always right before the closing\n// \"}\" paragraph: each new paragraph lands right above the anchor, below the\n// ones already inserted, which keeps the final order matching newParaTexts.\nfor (let i = 0; i < newParaTexts.length; i++) {\n  closingBrace.insertParagraph(newParaTexts[i], Word.InsertLocation.before);\n}\nawait context.sync();\n\n// --- 4. Append a closing \"```\" paragraph after the final \"]\" paragraph ---\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst allItems = paragraphs.items;\nconst lastPara = allItems[allItems.length - 1]; // the \"]\" paragraph\nlastPara.insertParagraph(\"```\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Wraps the JSON array of translation-review objects in a ```json fenced\n# code block, fixes the first entry's text, and appends two more entries.\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert a \"```json\" paragraph before the current first paragraph (\"[\") ---\n$d.Paragraphs(1).Range.InsertParagraphBefore()\n$d.Paragraphs(1).Range.Text = '```json'\n\n# Paragraph index map after the insertion above:\n# 1: ```json\n# 2: [\n# 3:   {\n# 4:     \"problematic_translated_sentence\": ...\n# 5:     \"respective_source_sentence\": ...\n# 6:     \"errors\": ...\n# 7:   }\n# 8: ]\n\n# --- 2. Fix the text of the first entry's three data paragraphs ---\n$d.Paragraphs(4).Range.Text = '    \"problematic_translated_sentence\": \"\u6240\u6709\u5728\u5f37\u8feb\u4e0b\u505a\u51fa\u7684\u52aa\u529b\u90fd\u9700\u8981\u72a7\u7272\u751f\u547d\u529b\u3002\",'\n$d.Paragraphs(5).Range.Text = '    \"respective_source_sentence\": \"Every effort under compulsion demands a sacrifice of life-energy.\",'\n$d.Paragraphs(6).Range.Text = '    \"errors\": \"Duplication of ''\u6240\u6709'' and ''\u505a\u51fa\u7684'' is redundant.\"'\n\n# --- 3. Insert the closing \"},\" of entry 1 plus two new entries, right\n#         before the lone \"  }\" paragraph that closes the JSON array ---\n$newParas = @(\n  '  },',\n  '  {',\n  '    \"problematic_translated_sentence\": \"\u6211\u5f9e\u672a\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9\u3002\",',\n  '    \"respective_source_sentence\": \"I never paid such a price.\",',\n  '    \"errors\": \"Mistranslation of ''paid such a price'' as ''\u4ed8\u51fa\u904e\u9019\u6a23\u7684\u4ee3\u50f9''.\"',\n  '  },',\n  '  {',\n  '    \"problematic_translated_sentence\": \"\u76f8\u53cd\uff0c\u6211\u5f9e\u6211\u7684\u601d\u7dd2\u4e2d\u7372\u76ca\u532a\u6dfa\u3002\",',\n  '    \"respective_source_sentence\": \"On the contrary, I have thrived on my thoughts.\",',\n  '    \"errors\": \"Mistranslation of ''thrived on my thoughts'' as ''\u5f9e\u6211\u7684\u601d\u7dd2\u4e2d\u7372\u76ca\u532a\u6dfa''.\"'\n)\n\n# The closing \"  }\" paragraph is currently paragraph 7 (index after step 2).\n$closeIndex = 7\nforeach ($txt in $newParas) {\n  $d.Paragraphs($closeIndex).Range.InsertParagraphBefore()\n  $d.Paragraphs($closeIndex).Range.Text = $txt\n  $closeIndex = $closeIndex + 1\n}\n\n# --- 4. Append a closing \"```\" paragraph after the final \"]\" paragraph ---\n$lastIndex = $d.Paragraphs.Count\n$d.Paragraphs($lastIndex).Range.InsertParagraphAfter()\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = '```'\n"}
